$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"
# Card number is digits-only but must stay text (matches source data, which
# stores it as a string, not a number) - lead with an apostrophe so Excel
# keeps it as text instead of parsing it as a numeric value.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 03.05.2024"

# Row 6
$ws.Range("B6").Value = "07.05."
$ws.Range("C6").Value = "08.05."
$ws.Range("D6").Value = "BURGER KING Ebermannstadt"
$ws.Range("E6").Value = "11,12-"

# Row 7
$ws.Range("B7").Value = "11.05."
$ws.Range("C7").Value = "12.05."
$ws.Range("D7").Value = "KARTENZ./11.05 EDEKA RO"
$ws.Range("E7").Value = "73,96-"

# Row 8
$ws.Range("B8").Value = "14.05."
$ws.Range("C8").Value = "15.05."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 50343289"
$ws.Range("E8").Value = "83,06-"

# Row 9
$ws.Range("B9").Value = "15.05."
$ws.Range("C9").Value = "16.05."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-14761906"
$ws.Range("E9").Value = "53,79-"

# Row 10
$ws.Range("B10").Value = "19.05."
$ws.Range("C10").Value = "20.05."
$ws.Range("D10").Value = "KARTENZ./19.05 REWE RO"
$ws.Range("E10").Value = "80,25-"

# Row 11 - previously an empty filler row, now gets a new transaction.
$ws.Range("B11").Value = "22.05."
$ws.Range("C11").Value = "23.05."
$ws.Range("D11").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E11").Value = "79,50-"
# E11 originally used the wrapped "empty amount" style (s=12); the filled-in
# amount cells above (E6:E10) use the narrower right-aligned style (s=17).
# Copy that formatting from the row above so E11 matches the others.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 25.05.2024"
$ws.Range("E12").Value = "381,68-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.06.2024"
